$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Chuckegg")
$ws.Range("B2").Value = 15
$ws.Range("C2").Value = -65
$ws.Range("E2").Value = -65
$ws.Range("G2").Value = -65
$ws.Range("I2").Value = -65
$ws.Range("B3").Value = 2505057.499999999
$ws.Range("C3").Value = 17489.16666666698
$ws.Range("E3").Value = 17489.16666666698
$ws.Range("G3").Value = 17489.16666666884
$ws.Range("I3").Value = 17489.16666666698
$ws.Range("B4").Value = 514
$ws.Range("C4").Value = 4
$ws.Range("E4").Value = 4
$ws.Range("G4").Value = 4
$ws.Range("I4").Value = 4
$ws.Range("B6").Value = 840387
$ws.Range("C6").Value = 4056
$ws.Range("E6").Value = 4056
$ws.Range("G6").Value = 4056
$ws.Range("I6").Value = 4056
$ws.Range("B7").Value = 9924
$ws.Range("C7").Value = 14
$ws.Range("E7").Value = 14
$ws.Range("G7").Value = 14
$ws.Range("I7").Value = 14
$ws.Range("B8").Value = 2714
$ws.Range("C8").Value = 3
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3
$ws.Range("I8").Value = 3
$ws.Range("B9").Value = 24705
$ws.Range("C9").Value = 125
$ws.Range("E9").Value = 125
$ws.Range("G9").Value = 125
$ws.Range("I9").Value = 125
$ws.Range("B10").Value = 6131
$ws.Range("C10").Value = 4
$ws.Range("E10").Value = 4
$ws.Range("G10").Value = 4
$ws.Range("I10").Value = 4
$ws.Range("B11").Value = 309039
$ws.Range("C11").Value = 1371
$ws.Range("E11").Value = 1371
$ws.Range("G11").Value = 1371
$ws.Range("I11").Value = 1371
$ws.Range("B13").Value = 6818
$ws.Range("C13").Value = 11
$ws.Range("E13").Value = 11
$ws.Range("G13").Value = 11
$ws.Range("I13").Value = 11
$ws.Range("B14").Value = 18556
$ws.Range("C14").Value = 121
$ws.Range("E14").Value = 121
$ws.Range("G14").Value = 121
$ws.Range("I14").Value = 121
$ws.Range("B15").Value = 36053
$ws.Range("C15").Value = 211
$ws.Range("E15").Value = 211
$ws.Range("G15").Value = 211
$ws.Range("I15").Value = 211
$ws.Range("B16").Value = 22943
$ws.Range("C16").Value = 128
$ws.Range("E16").Value = 128
$ws.Range("G16").Value = 128
$ws.Range("I16").Value = 128
$ws.Range("B17").Value = 17809
$ws.Range("C17").Value = 98
$ws.Range("E17").Value = 98
$ws.Range("G17").Value = 98
$ws.Range("I17").Value = 98
$ws.Range("B18").Value = 11585
$ws.Range("C18").Value = 71
$ws.Range("E18").Value = 71
$ws.Range("G18").Value = 71
$ws.Range("I18").Value = 71
$ws.Range("B21").Value = 1208
$ws.Range("C21").Value = 9
$ws.Range("E21").Value = 9
$ws.Range("G21").Value = 9
$ws.Range("I21").Value = 9
$ws.Range("B22").Value = 2846128
$ws.Range("C22").Value = 14170
$ws.Range("E22").Value = 14170
$ws.Range("G22").Value = 14170
$ws.Range("I22").Value = 14170
$ws = $wb.Worksheets.Item("Kimmiii20")
$ws.Range("E2").Value = 60
$ws.Range("G2").Value = 60
$ws.Range("I2").Value = 60
$ws.Range("B2").Value = 65
$ws.Range("C2").Value = 60
$ws.Range("B3").Value = 189387.9166666667
$ws.Range("C3").Value = 37112.50000000017
$ws.Range("E3").Value = 37112.50000000017
$ws.Range("G3").Value = 37112.50000000017
$ws.Range("I3").Value = 37112.50000000017
$ws.Range("B6").Value = 39546
$ws.Range("C6").Value = 8186
$ws.Range("E6").Value = 8186
$ws.Range("G6").Value = 8186
$ws.Range("I6").Value = 8186
$ws.Range("B7").Value = 1108
$ws.Range("C7").Value = 158
$ws.Range("E7").Value = 158
$ws.Range("G7").Value = 158
$ws.Range("I7").Value = 158
$ws.Range("B8").Value = 177
$ws.Range("C8").Value = 26
$ws.Range("E8").Value = 26
$ws.Range("G8").Value = 26
$ws.Range("I8").Value = 26
$ws.Range("B9").Value = 2126
$ws.Range("C9").Value = 460
$ws.Range("E9").Value = 460
$ws.Range("G9").Value = 460
$ws.Range("I9").Value = 460
$ws.Range("B10").Value = 789
$ws.Range("C10").Value = 97
$ws.Range("E10").Value = 97
$ws.Range("G10").Value = 97
$ws.Range("I10").Value = 97
$ws.Range("B11").Value = 27112
$ws.Range("C11").Value = 5010
$ws.Range("E11").Value = 5010
$ws.Range("G11").Value = 5010
$ws.Range("I11").Value = 5010
$ws.Range("B13").Value = 775
$ws.Range("C13").Value = 124
$ws.Range("E13").Value = 124
$ws.Range("G13").Value = 124
$ws.Range("I13").Value = 124
$ws.Range("B14").Value = 1337
$ws.Range("C14").Value = 363
$ws.Range("E14").Value = 363
$ws.Range("G14").Value = 363
$ws.Range("I14").Value = 363
$ws.Range("B15").Value = 1834
$ws.Range("C15").Value = 326
$ws.Range("E15").Value = 326
$ws.Range("G15").Value = 326
$ws.Range("I15").Value = 326
$ws.Range("B16").Value = 1167
$ws.Range("C16").Value = 210
$ws.Range("E16").Value = 210
$ws.Range("G16").Value = 210
$ws.Range("I16").Value = 210
$ws.Range("B17").Value = 599
$ws.Range("C17").Value = 148
$ws.Range("E17").Value = 148
$ws.Range("G17").Value = 148
$ws.Range("I17").Value = 148
$ws.Range("B18").Value = 458
$ws.Range("C18").Value = 110
$ws.Range("E18").Value = 110
$ws.Range("G18").Value = 110
$ws.Range("I18").Value = 110
$ws.Range("B22").Value = 223827
$ws.Range("C22").Value = 48546
$ws.Range("E22").Value = 48546
$ws.Range("G22").Value = 48546
$ws.Range("I22").Value = 48546
$ws = $wb.Worksheets.Item("I_HAMSTER")
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 4
$ws.Range("E2").Value = 4
$ws.Range("G2").Value = 4
$ws.Range("I2").Value = 4
$ws.Range("B3").Value = 213425.8333333331
$ws.Range("C3").Value = 31735
$ws.Range("E3").Value = 31735
$ws.Range("G3").Value = 31735
$ws.Range("I3").Value = 31735
$ws.Range("B4").Value = 45
$ws.Range("C4").Value = 6
$ws.Range("E4").Value = 6
$ws.Range("G4").Value = 6
$ws.Range("I4").Value = 6
$ws.Range("B5").Value = 59760
$ws.Range("C5").Value = 3520
$ws.Range("E5").Value = 3520
$ws.Range("G5").Value = 3520
$ws.Range("I5").Value = 3520
$ws.Range("B6").Value = 75006
$ws.Range("C6").Value = 9377
$ws.Range("E6").Value = 9377
$ws.Range("G6").Value = 9377
$ws.Range("I6").Value = 9377
$ws.Range("B7").Value = 1116
$ws.Range("C7").Value = 107
$ws.Range("E7").Value = 107
$ws.Range("G7").Value = 107
$ws.Range("I7").Value = 107
$ws.Range("B8").Value = 142
$ws.Range("C8").Value = 21
$ws.Range("E8").Value = 21
$ws.Range("G8").Value = 21
$ws.Range("I8").Value = 21
$ws.Range("B9").Value = 2418
$ws.Range("C9").Value = 323
$ws.Range("E9").Value = 323
$ws.Range("G9").Value = 323
$ws.Range("I9").Value = 323
$ws.Range("B10").Value = 648
$ws.Range("C10").Value = 64
$ws.Range("E10").Value = 64
$ws.Range("G10").Value = 64
$ws.Range("I10").Value = 64
$ws.Range("B11").Value = 30195
$ws.Range("C11").Value = 2988
$ws.Range("E11").Value = 2988
$ws.Range("G11").Value = 2988
$ws.Range("I11").Value = 2988
$ws.Range("B12").Value = 6
$ws.Range("C12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("I12").Value = 1
$ws.Range("B13").Value = 405
$ws.Range("C13").Value = 82
$ws.Range("E13").Value = 82
$ws.Range("G13").Value = 82
$ws.Range("I13").Value = 82
$ws.Range("B14").Value = 1770
$ws.Range("C14").Value = 259
$ws.Range("E14").Value = 259
$ws.Range("G14").Value = 259
$ws.Range("I14").Value = 259
$ws.Range("B15").Value = 3918
$ws.Range("C15").Value = 418
$ws.Range("E15").Value = 418
$ws.Range("G15").Value = 418
$ws.Range("I15").Value = 418
$ws.Range("B16").Value = 1204
$ws.Range("C16").Value = 246
$ws.Range("E16").Value = 246
$ws.Range("G16").Value = 246
$ws.Range("I16").Value = 246
$ws.Range("B17").Value = 765
$ws.Range("C17").Value = 128
$ws.Range("E17").Value = 128
$ws.Range("G17").Value = 128
$ws.Range("I17").Value = 128
$ws.Range("B18").Value = 541
$ws.Range("C18").Value = 131
$ws.Range("E18").Value = 131
$ws.Range("G18").Value = 131
$ws.Range("I18").Value = 131
$ws.Range("B19").Value = 102
$ws.Range("C19").Value = 19
$ws.Range("E19").Value = 19
$ws.Range("G19").Value = 19
$ws.Range("I19").Value = 19
$ws.Range("B21").Value = 98
$ws.Range("C21").Value = 22
$ws.Range("E21").Value = 22
$ws.Range("G21").Value = 22
$ws.Range("I21").Value = 22
$ws.Range("B22").Value = 158828
$ws.Range("C22").Value = 36060
$ws.Range("E22").Value = 36060
$ws.Range("G22").Value = 36060
$ws.Range("I22").Value = 36060
$ws = $wb.Worksheets.Item("one4kat")
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = -50
$ws.Range("E2").Value = -2
$ws.Range("G2").Value = -9
$ws.Range("I2").Value = -50
$ws.Range("B5").Value = 142545
$ws.Range("C5").Value = 1233
$ws.Range("E5").Value = 22
$ws.Range("G5").Value = 3545
$ws.Range("I5").Value = 1233
$ws = $wb.Worksheets.Item("pogdvde")
$ws.Range("B3").Value = 1673969.583333331
$ws.Range("C3").Value = 33453.74999999814
$ws.Range("E3").Value = 33453.74999999814
$ws.Range("G3").Value = 33453.74999999814
$ws.Range("I3").Value = 33453.74999999814
$ws.Range("B5").Value = 174837
$ws.Range("C5").Value = 941
$ws.Range("E5").Value = 941
$ws.Range("G5").Value = 941
$ws.Range("I5").Value = 941
$ws.Range("B6").Value = 627196
$ws.Range("C6").Value = 13481
$ws.Range("E6").Value = 13481
$ws.Range("G6").Value = 13481
$ws.Range("I6").Value = 13481
$ws.Range("B7").Value = 8953
$ws.Range("C7").Value = 90
$ws.Range("E7").Value = 90
$ws.Range("G7").Value = 90
$ws.Range("I7").Value = 90
$ws.Range("B8").Value = 2232
$ws.Range("C8").Value = 20
$ws.Range("E8").Value = 20
$ws.Range("G8").Value = 20
$ws.Range("I8").Value = 20
$ws.Range("B9").Value = 19514
$ws.Range("C9").Value = 415
$ws.Range("E9").Value = 415
$ws.Range("G9").Value = 415
$ws.Range("I9").Value = 415
$ws.Range("B10").Value = 6227
$ws.Range("C10").Value = 55
$ws.Range("E10").Value = 55
$ws.Range("G10").Value = 55
$ws.Range("I10").Value = 55
$ws.Range("B11").Value = 241406
$ws.Range("C11").Value = 4776
$ws.Range("E11").Value = 4776
$ws.Range("G11").Value = 4776
$ws.Range("I11").Value = 4776
$ws.Range("B13").Value = 6380
$ws.Range("C13").Value = 68
$ws.Range("E13").Value = 68
$ws.Range("G13").Value = 68
$ws.Range("I13").Value = 68
$ws.Range("B14").Value = 13262
$ws.Range("C14").Value = 359
$ws.Range("E14").Value = 359
$ws.Range("G14").Value = 359
$ws.Range("I14").Value = 359
$ws.Range("B15").Value = 29305
$ws.Range("C15").Value = 664
$ws.Range("E15").Value = 664
$ws.Range("G15").Value = 664
$ws.Range("I15").Value = 664
$ws.Range("B16").Value = 20177
$ws.Range("C16").Value = 441
$ws.Range("E16").Value = 441
$ws.Range("G16").Value = 441
$ws.Range("I16").Value = 441
$ws.Range("B21").Value = 1049
$ws.Range("C21").Value = 23
$ws.Range("E21").Value = 23
$ws.Range("G21").Value = 23
$ws.Range("I21").Value = 23
$ws.Range("B22").Value = 1944513
$ws.Range("C22").Value = 40581
$ws.Range("E22").Value = 40581
$ws.Range("G22").Value = 40581
$ws.Range("I22").Value = 40581
$ws = $wb.Worksheets.Item("ewz_")
$ws.Range("B2").Value = 17
$ws.Range("C2").Value = 12
$ws.Range("E2").Value = 12
$ws.Range("G2").Value = 12
$ws.Range("I2").Value = 12
$ws.Range("B3").Value = 20127.49999999999
$ws.Range("C3").Value = 335.8333333333285
$ws.Range("E3").Value = 335.8333333333285
$ws.Range("G3").Value = 335.8333333333285
$ws.Range("I3").Value = 335.8333333333285
$ws.Range("B4").Value = 7
$ws.Range("C4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("B6").Value = 3214
$ws.Range("C6").Value = 83
$ws.Range("E6").Value = 83
$ws.Range("G6").Value = 83
$ws.Range("I6").Value = 83
$ws.Range("B7").Value = 199
$ws.Range("C7").Value = 4
$ws.Range("E7").Value = 4
$ws.Range("G7").Value = 4
$ws.Range("I7").Value = 4
$ws.Range("B8").Value = 37
$ws.Range("C8").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("B9").Value = 288
$ws.Range("C9").Value = 6
$ws.Range("E9").Value = 6
$ws.Range("G9").Value = 6
$ws.Range("I9").Value = 6
$ws.Range("B10").Value = 104
$ws.Range("C10").Value = 4
$ws.Range("E10").Value = 4
$ws.Range("G10").Value = 4
$ws.Range("I10").Value = 4
$ws.Range("B11").Value = 2437
$ws.Range("C11").Value = 48
$ws.Range("E11").Value = 48
$ws.Range("G11").Value = 48
$ws.Range("I11").Value = 48
$ws.Range("B12").Value = 154
$ws.Range("C12").Value = 3
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 3
$ws.Range("I12").Value = 3
$ws.Range("B13").Value = 184
$ws.Range("C13").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 2
$ws.Range("I13").Value = 2
$ws.Range("B14").Value = 133
$ws.Range("C14").Value = 4
$ws.Range("E14").Value = 4
$ws.Range("G14").Value = 4
$ws.Range("I14").Value = 4
$ws.Range("B15").Value = 97
$ws.Range("C15").Value = 3
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3
$ws.Range("I15").Value = 3
$ws.Range("B16").Value = 48
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("G16").Value = 1
$ws.Range("I16").Value = 1
$ws.Range("B19").Value = 3
$ws.Range("C19").Value = 1
$ws.Range("E19").Value = 1
$ws.Range("G19").Value = 1
$ws.Range("I19").Value = 1
$ws.Range("B20").Value = 31074
$ws.Range("C20").Value = 566
$ws.Range("E20").Value = 566
$ws.Range("G20").Value = 566
$ws.Range("I20").Value = 566
